# Update "想去人数" (number interested) counts on the 展览 and 全部类型 sheets
# to reflect freshly generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 4559
$ws1.Range("F5").Value = 42
$ws1.Range("F6").Value = 468

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 4559
$ws4.Range("F7").Value = 42
$ws4.Range("F8").Value = 468
